# Weekly update: a new price record ("Poroto granado", Macroferia Regional de
# Talca) is inserted as row 78, pushing the existing rows 78-113 down to
# 79-114 (each keeps its own Fecha/Volumen/Precio values, just shifted one
# row down). The sheet's used range grows from A1:R113 to A1:R114.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank row above the current row 78 - everything below
# (rows 78..113) shifts down to 79..114.
$ws.Rows.Item(78).Insert()

# Populate the newly inserted row 78 with the new weekly record.
$ws.Range("A78").Value2 = 5
$ws.Range("B78").Value2 = "Macroferia Regional de Talca"
$ws.Range("C78").Value2 = "Maule"
$ws.Range("D78").Value2 = 44609
$ws.Range("E78").Value2 = 7
$ws.Range("F78").Value2 = 100112030
$ws.Range("G78").Value2 = "Poroto granado"
$ws.Range("H78").Value2 = "Sin especificar"
$ws.Range("I78").Value2 = "Primera"
$ws.Range("J78").Value2 = 400
$ws.Range("K78").Value2 = 17000
$ws.Range("L78").Value2 = 17000
$ws.Range("M78").Value2 = 17000
$ws.Range("N78").Value2 = "$/saco 25 kilos"
$ws.Range("O78").Value2 = "Región del Maule"
$ws.Range("P78").Value2 = 680
$ws.Range("Q78").Value2 = 25
$ws.Range("R78").Value2 = "Hortaliza"
